$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1089330235413494
$ws.Range("H2").Value = 15.00423028218035
$ws.Range("I2").Value = 70.27997336718379
$ws.Range("G3").Value = 0.08633469868054416
$ws.Range("H3").Value = 30.62429424758897
$ws.Range("G4").Value = 0.05564443211082127
$ws.Range("H4").Value = 493.5236702004174
$ws.Range("G5").Value = 0.06473273844841052
$ws.Range("H5").Value = 583.4476899099353
$ws.Range("G6").Value = -0.2480781787696523
$ws.Range("H6").Value = -12.15326166811768
$ws.Range("G7").Value = -0.1938512765520452
$ws.Range("H7").Value = 22.4209152946375
$ws.Range("G8").Value = -0.348482098000785
$ws.Range("H8").Value = 5.880736411723247
$ws.Range("G9").Value = -0.4465935199590182
$ws.Range("H9").Value = -12.03086184648177
$ws.Range("G10").Value = -0.03661058812501135
$ws.Range("H10").Value = -325.9231135764834
$ws.Range("G11").Value = 0.1109858677711741
$ws.Range("H11").Value = 790.8177784076412
$ws.Range("G12").Value = 0.2442780973087719
$ws.Range("H12").Value = 7.522068059452092
$ws.Range("G13").Value = 0.3048948402577875
$ws.Range("H13").Value = 15.77966240322421
$ws.Range("G14").Value = 0.003377446125900369
$ws.Range("H14").Value = 135.3020797141205
$ws.Range("G15").Value = 0.01807478892439136
$ws.Range("H15").Value = -10.46196874438527
$ws.Range("G16").Value = 0.1471397429223727
$ws.Range("H16").Value = 24.69564278834616
$ws.Range("G17").Value = 0.2019013891914626
$ws.Range("H17").Value = -7.737196508946391
$ws.Range("G18").Value = 0.05066496190056213
$ws.Range("H18").Value = -16.21000852522241
$ws.Range("G19").Value = 0.06335529464375943
$ws.Range("H19").Value = -29.67500940190437
$ws.Range("G20").Value = -0.1493376391708737
$ws.Range("H20").Value = -2.615707229767667
$ws.Range("G21").Value = -0.1655791104804148
$ws.Range("H21").Value = 17.13996972952506
$ws.Range("G22").Value = 0.06220922636670745
$ws.Range("H22").Value = 14.37991635058684
$ws.Range("G23").Value = 0.04744535075041098
$ws.Range("H23").Value = 16.17331556361706
$ws.Range("G24").Value = 0.1240932540643074
$ws.Range("H24").Value = 7.222516877577907
$ws.Range("G25").Value = 0.1450114575138241
$ws.Range("H25").Value = -4.64542053481577
$ws.Range("G26").Value = 0.03522127456327665
$ws.Range("H26").Value = -33.38991132908583
$ws.Range("G27").Value = 0.0270347553693973
$ws.Range("H27").Value = -46.43062100028067
$ws.Range("G28").Value = 0.1431970712570723
$ws.Range("H28").Value = -6.352712753401884
$ws.Range("G29").Value = 0.1563239313775818
$ws.Range("H29").Value = -8.42348942306189
$ws.Range("G30").Value = 0.04406200190503957
$ws.Range("H30").Value = 125.1957797435301
$ws.Range("G31").Value = 0.01911163817547751
$ws.Range("H31").Value = 96.92628611492408
$ws.Range("G32").Value = 0.007069033143765936
$ws.Range("H32").Value = -81.0447496072093
$ws.Range("G33").Value = 0.02113011879896545
$ws.Range("H33").Value = -19.05297608858366
$ws.Range("G34").Value = 0.1135768874065771
$ws.Range("H34").Value = -11.24587416545573
$ws.Range("G35").Value = 0.1247405251258358
$ws.Range("H35").Value = -3.046828740105326
$ws.Range("G36").Value = -0.05648232002048782
$ws.Range("H36").Value = -475.7500191276877
$ws.Range("G37").Value = 0.0176297342143285
$ws.Range("H37").Value = 15.11800949944664
$ws.Range("G38").Value = -0.01247581372136387
$ws.Range("H38").Value = -510.8367958859019
$ws.Range("G39").Value = -0.03503784162529756
$ws.Range("H39").Value = -4.874790846011702
$ws.Range("G40").Value = 0.1272225312727454
$ws.Range("H40").Value = -13.77607567012815
$ws.Range("G41").Value = 0.1558683796813407
$ws.Range("H41").Value = -3.426799027935116
$ws.Range("G42").Value = 0.03960050263682657
$ws.Range("H42").Value = -38.66559065016072
$ws.Range("G43").Value = 0.07230136851666458
$ws.Range("H43").Value = 107.9985048403623
$ws.Range("G44").Value = 0.03065116904809738
$ws.Range("H44").Value = 117.1872816745017
$ws.Range("G45").Value = 0.02291328482669539
$ws.Range("H45").Value = -44.19293840450187
$ws.Range("G46").Value = -0.03469376567376054
$ws.Range("H46").Value = 47.28982474483102
$ws.Range("G47").Value = -0.02288662598077324
$ws.Range("H47").Value = 44.5974567142873
$ws.Range("G48").Value = -0.1208150834097541
$ws.Range("H48").Value = 4.096296282550468
$ws.Range("G49").Value = -0.115229025422366
$ws.Range("H49").Value = 41.65069643040344
$ws.Range("G50").Value = 0.08714110998930936
$ws.Range("H50").Value = -19.96124703255114
$ws.Range("G51").Value = 0.1134302532499808
$ws.Range("H51").Value = 13.12388172783215
$ws.Range("G52").Value = 0.06152233664544582
$ws.Range("H52").Value = 3.192407917037723
$ws.Range("G53").Value = 0.0826345111258542
$ws.Range("H53").Value = 22.34413941261884
$ws.Range("G54").Value = -0.06681443655767105
$ws.Range("H54").Value = 4.442388431988669
$ws.Range("G55").Value = -0.07234877222006822
$ws.Range("H55").Value = 6.327900893646269
$ws.Range("G56").Value = 0.09498575403203936
$ws.Range("H56").Value = 107.2744762529767
$ws.Range("G57").Value = 0.1468825598540809
$ws.Range("H57").Value = 2740.902659428968

Write-Host "done"